# Auto-generated edit script: updates cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.253.37"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "3.493.86"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'587.13"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "'134.56"
$ws.Range("E6").Value = "  +1.91%  "
$ws.Range("D7").Value = "3.494.26"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "4.087.90"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "3.495.09"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "64.313.47"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "'25.43"
$ws.Range("E18").Value = "  -8.60%  "
$ws.Range("D19").Value = "'9.87"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "'13.61"
$ws.Range("E21").Value = "  -6.24%  "
$ws.Range("D22").Value = "'388.61"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").Value = "'0.566"
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("D24").Value = "3.634.14"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "'74.30"
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").Value = "'0.0000114"
$ws.Range("E28").Value = "  +1.30%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("E31").Value = "  -5.05%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").Value = "'8.22"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").Value = "3.516.36"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("E36").Value = "  +3.63%  "
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("D41").Value = "'162.50"
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'25.44"
$ws.Range("E45").Value = "  -5.63%  "
$ws.Range("D46").Value = "'41.81"
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").Value = "'4.41"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").Value = "2.472.10"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("E51").Value = "  -2.40%  "
